$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efnb2"
$ws.Cells.Item(2,3).Value = "Rhbdl2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 37.98277566666666
$ws.Cells.Item(2,8).Value = 113.948327
$ws.Cells.Item(2,9).Value = 0.697850645410475
$ws.Cells.Item(2,10).Value = 0.6978506454104751
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.082435
$ws.Cells.Item(2,14).Value = 6.247305
$ws.Cells.Item(2,15).Value = 0.9920259111440977
$ws.Cells.Item(2,16).Value = 0.9920259111440977
$ws.Cells.Item(2,17).Value = 79.09666144541498
$ws.Cells.Item(2,18).Value = 711.8699530087349
$ws.Cells.Item(2,19).Value = 0.6922859223558231
$ws.Cells.Item(2,20).Value = 0.6922859223558232

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efnb2"
$ws.Cells.Item(3,3).Value = "Rhbdl2"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 37.98277566666666
$ws.Cells.Item(3,8).Value = 113.948327
$ws.Cells.Item(3,9).Value = 0.697850645410475
$ws.Cells.Item(3,10).Value = 0.6978506454104751
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.016739
$ws.Cells.Item(3,14).Value = 0.050217
$ws.Cells.Item(3,15).Value = 0.00797408885590237
$ws.Cells.Item(3,16).Value = 0.007974088855902369
$ws.Cells.Item(3,17).Value = 0.6357936818843333
$ws.Cells.Item(3,18).Value = 5.722143136959
$ws.Cells.Item(3,19).Value = 0.005564723054651946
$ws.Cells.Item(3,20).Value = 0.005564723054651945

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Efnb2"
$ws.Cells.Item(4,3).Value = "Rhbdl2"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 9.685730333333334
$ws.Cells.Item(4,8).Value = 29.057191
$ws.Cells.Item(4,9).Value = 0.1779541659542351
$ws.Cells.Item(4,10).Value = 0.1779541659542352
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.082435
$ws.Cells.Item(4,14).Value = 6.247305
$ws.Cells.Item(4,15).Value = 0.9920259111440977
$ws.Cells.Item(4,16).Value = 0.9920259111440977
$ws.Cells.Item(4,17).Value = 20.169903846695
$ws.Cells.Item(4,18).Value = 181.529134620255
$ws.Cells.Item(4,19).Value = 0.1765351436226381
$ws.Cells.Item(4,20).Value = 0.1765351436226381

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Efnb2"
$ws.Cells.Item(5,3).Value = "Rhbdl2"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 9.685730333333334
$ws.Cells.Item(5,8).Value = 29.057191
$ws.Cells.Item(5,9).Value = 0.1779541659542351
$ws.Cells.Item(5,10).Value = 0.1779541659542352
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.016739
$ws.Cells.Item(5,14).Value = 0.050217
$ws.Cells.Item(5,15).Value = 0.00797408885590237
$ws.Cells.Item(5,16).Value = 0.007974088855902369
$ws.Cells.Item(5,17).Value = 0.1621294400496667
$ws.Cells.Item(5,18).Value = 1.459164960447
$ws.Cells.Item(5,19).Value = 0.001419022331597067
$ws.Cells.Item(5,20).Value = 0.001419022331597067

# Row 6
$ws.Cells.Item(6,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value = "Efnb2"
$ws.Cells.Item(6,3).Value = "Rhbdl2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.5676613333333332
$ws.Cells.Item(6,8).Value = 1.702984
$ws.Cells.Item(6,9).Value = 0.01042953867610283
$ws.Cells.Item(6,10).Value = 0.01042953867610283
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.082435
$ws.Cells.Item(6,14).Value = 6.247305
$ws.Cells.Item(6,15).Value = 0.9920259111440977
$ws.Cells.Item(6,16).Value = 0.9920259111440977
$ws.Cells.Item(6,17).Value = 1.18211782868
$ws.Cells.Item(6,18).Value = 10.63906045812
$ws.Cells.Item(6,19).Value = 0.01034637260797352
$ws.Cells.Item(6,20).Value = 0.01034637260797352

# Row 7
$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Efnb2"
$ws.Cells.Item(7,3).Value = "Rhbdl2"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.5676613333333332
$ws.Cells.Item(7,8).Value = 1.702984
$ws.Cells.Item(7,9).Value = 0.01042953867610283
$ws.Cells.Item(7,10).Value = 0.01042953867610283
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.016739
$ws.Cells.Item(7,14).Value = 0.050217
$ws.Cells.Item(7,15).Value = 0.00797408885590237
$ws.Cells.Item(7,16).Value = 0.007974088855902369
$ws.Cells.Item(7,17).Value = 0.009502083058666665
$ws.Cells.Item(7,18).Value = 0.08551874752799998
$ws.Cells.Item(7,19).Value = 0.00008316606812931435
$ws.Cells.Item(7,20).Value = 0.00008316606812931433

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Efnb2"
$ws.Cells.Item(8,3).Value = "Rhbdl2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 5.823095333333334
$ws.Cells.Item(8,8).Value = 17.469286
$ws.Cells.Item(8,9).Value = 0.1069866739681064
$ws.Cells.Item(8,10).Value = 0.1069866739681064
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.082435
$ws.Cells.Item(8,14).Value = 6.247305
$ws.Cells.Item(8,15).Value = 0.9920259111440977
$ws.Cells.Item(8,16).Value = 0.9920259111440977
$ws.Cells.Item(8,17).Value = 12.12621753047
$ws.Cells.Item(8,18).Value = 109.13595777423
$ws.Cells.Item(8,19).Value = 0.1061335527234873
$ws.Cells.Item(8,20).Value = 0.1061335527234873

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Efnb2"
$ws.Cells.Item(9,3).Value = "Rhbdl2"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 5.823095333333334
$ws.Cells.Item(9,8).Value = 17.469286
$ws.Cells.Item(9,9).Value = 0.1069866739681064
$ws.Cells.Item(9,10).Value = 0.1069866739681064
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.016739
$ws.Cells.Item(9,14).Value = 0.050217
$ws.Cells.Item(9,15).Value = 0.00797408885590237
$ws.Cells.Item(9,16).Value = 0.007974088855902369
$ws.Cells.Item(9,17).Value = 0.09747279278466668
$ws.Cells.Item(9,18).Value = 0.8772551350619999
$ws.Cells.Item(9,19).Value = 0.0008531212446191378
$ws.Cells.Item(9,20).Value = 0.0008531212446191375

# Row 10
$ws.Cells.Item(10,1).Value = "Resolving-Mac"
$ws.Cells.Item(10,2).Value = "Efnb2"
$ws.Cells.Item(10,3).Value = "Rhbdl2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.3689676666666666
$ws.Cells.Item(10,8).Value = 1.106903
$ws.Cells.Item(10,9).Value = 0.006778975991080511
$ws.Cells.Item(10,10).Value = 0.006778975991080512
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.082435
$ws.Cells.Item(10,14).Value = 6.247305
$ws.Cells.Item(10,15).Value = 0.9920259111440977
$ws.Cells.Item(10,16).Value = 0.9920259111440977
$ws.Cells.Item(10,17).Value = 0.7683511829349998
$ws.Cells.Item(10,18).Value = 6.915160646415
$ws.Cells.Item(10,19).Value = 0.006724919834175606
$ws.Cells.Item(10,20).Value = 0.006724919834175607

# Row 11
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Efnb2"
$ws.Cells.Item(11,3).Value = "Rhbdl2"
$ws.Cells.Item(11,4).Value = "MuSCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.3689676666666666
$ws.Cells.Item(11,8).Value = 1.106903
$ws.Cells.Item(11,9).Value = 0.006778975991080511
$ws.Cells.Item(11,10).Value = 0.006778975991080512
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.016739
$ws.Cells.Item(11,14).Value = 0.050217
$ws.Cells.Item(11,15).Value = 0.00797408885590237
$ws.Cells.Item(11,16).Value = 0.007974088855902369
$ws.Cells.Item(11,17).Value = 0.006176149772333333
$ws.Cells.Item(11,18).Value = 0.05558534795099999
$ws.Cells.Item(11,19).Value = 0.00005405615690490483
$ws.Cells.Item(11,20).Value = 0.00005405615690490483
